$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "['ifs-fesom', 'icon', 'ifs-nemo-er', 'hadgem3-mediumres', 'hadgem3-hires', 'hadgem3-lowres']"
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,18,19,20,21)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = $newValue
}
